$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "d1"
$ws.Range("A2").Value = "d2"
$ws.Range("A3").Value = "d3"
$ws.Range("A4").Value = "d4"

$ws.Range("A4").Select()
